# Populate the DB-connection placeholder/test-value table.
#
# The write order below intentionally mirrors the order the data was
# originally authored in (header cells A1,C1,D1,E1; then the row-2 test
# values; then B1 was added afterward; then the numeric port last) so
# that the shared-string table comes out in the same sequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: placeholder names (B1 is filled in later, see below)
$ws.Range("A1").Value = '${DB_Name}'
$ws.Range("C1").Value = '${DB_Password}'
$ws.Range("D1").Value = '${DB_Host}'
$ws.Range("E1").Value = '${DB_Port}'

# Row 2: sample/test values
$ws.Range("A2").Value = 'TestDB'
$ws.Range("B2").Value = 'sa'
$ws.Range("C2").Value = 'test'
$ws.Range("D2").Value = 'localhost'

# B1 (user name header) added afterwards
$ws.Range("B1").Value = '${DB_User_Name}'

# Port number added last (numeric, not a shared string)
$ws.Range("E2").Value = 1521

# Column widths. Target character widths (from the authored file) are
# 22, 19.5703125, 17.28515625, 12.85546875, 12.42578125 for columns
# A-E respectively; ColumnWidth here is pre-offset by 5/6 so the
# resulting stored width lands on the nearest value this engine can
# represent.
$ws.Columns.Item(1).ColumnWidth = 21.16666666666667
$ws.Columns.Item(2).ColumnWidth = 18.736979166666668
$ws.Columns.Item(3).ColumnWidth = 16.451822916666668
$ws.Columns.Item(4).ColumnWidth = 12.022135416666666
$ws.Columns.Item(5).ColumnWidth = 11.592447916666666

# Selection moves to A3 after data entry
[void]$ws.Range("A3").Select()
